$wb = $excel.ActiveWorkbook
$linkedIn = $wb.Worksheets.Item("LinkedIn")

# Add "Wellfound" as a copy of "LinkedIn" placed right after it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$linkedIn.Copy($null, $lastSheet)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "Wellfound"

# Add "Others" as a copy of "LinkedIn" placed at the end.
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$linkedIn.Copy($null, $lastSheet2)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "Others"
